# Insert a new "Team Scalability" slide before the current slide 13
# ("Maintainability at Scale"), pushing it and all the later slides
# down by one position. The new slide reuses the same Blank layout
# (layout index 7) and the same shape layout (title textbox, accent
# rectangle bar, bullet textbox) used by the surrounding slides.

$p = $ppt.ActivePresentation

# Insert the new slide at position 13, using the Blank layout (the
# same one used by every other content slide in this deck).
$s = $p.Slides.Add(13, 7)

# --- Title textbox -------------------------------------------------
$title = $s.Shapes.AddTextbox(1, 36, 28.8, 648, 57.6)
$title.TextFrame.TextRange.Text = "Team Scalability"
$title.TextFrame.TextRange.Font.Size = 32
$title.TextFrame.TextRange.Font.Bold = $true
$title.TextFrame.TextRange.Font.Color.RGB = 3021338
$title.TextFrame.WordWrap = $false
$title.TextFrame.AutoSize = 1
$title.Left = 36
$title.Top = 28.8
$title.Width = 648
$title.Height = 57.6

# --- Accent bar rectangle -------------------------------------------
$rect = $s.Shapes.AddShape(1, 36, 82.8, 648, 2.16)
$rect.Fill.ForeColor.RGB = 13400576
$rect.Line.Visible = $false

# --- Bullet body textbox --------------------------------------------
$body = $s.Shapes.AddTextbox(1, 36, 100.8, 648, 324)
$body.TextFrame.TextRange.Text = "• Work can be parallelized across team members`r• Endpoint development decoupled from Ansible module work`r• Clear module boundaries reduce merge conflicts`r• New contributors can focus on one domain without full codebase knowledge`r• Enables efficient resource allocation across projects"
$body.TextFrame.TextRange.Font.Size = 20
$body.TextFrame.TextRange.Font.Color.RGB = 3355443
$body.TextFrame.TextRange.ParagraphFormat.SpaceAfter = 12
$body.TextFrame.WordWrap = $true
$body.TextFrame.AutoSize = 1
$body.Left = 36
$body.Top = 100.8
$body.Width = 648
$body.Height = 324
